$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) stores plain numeric-looking text (e.g. "217.78",
# "26.228.72"). Writing such strings straight into .Value lets Excel
# auto-detect them as numbers (dropping trailing zeros, changing the
# stored type away from text). Force the whole data range to Text format
# first so every write lands as a literal string, then clear the number
# format again afterwards so the cells keep their original (unstyled)
# appearance.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.228.72"
$ws.Range("E2").Value = "  -1.24%  "
$ws.Range("D3").Value = "1.672.77"
$ws.Range("E3").Value = "  +1.24%  "
$ws.Range("E4").Value = "  -0.35%  "
$ws.Range("D5").Value = "217.78"
$ws.Range("E5").Value = "  -0.92%  "
$ws.Range("D6").Value = "0.5138"
$ws.Range("E6").Value = "  +1.49%  "
$ws.Range("D7").Value = "1.006"
$ws.Range("E7").Value = "  -0.33%  "
$ws.Range("D8").Value = "0.2661"
$ws.Range("E8").Value = "  +5.07%  "
$ws.Range("D9").Value = "0.06389"
$ws.Range("E9").Value = "  +4.10%  "
$ws.Range("D10").Value = "21.60"
$ws.Range("E10").Value = "  +0.43%  "
$ws.Range("D11").Value = "0.07392"
$ws.Range("E11").Value = "  +0.52%  "
$ws.Range("D12").Value = "1.676.60"
$ws.Range("E12").Value = "  +1.37%  "
$ws.Range("E13").Value = "  +3.15%  "
$ws.Range("D14").Value = "0.5834"
$ws.Range("E14").Value = "  +2.04%  "
$ws.Range("D15").Value = "1.900.80"
$ws.Range("E15").Value = "  +1.28%  "
$ws.Range("D16").Value = "0.000008697"
$ws.Range("E16").Value = "  +8.63%  "
$ws.Range("D17").Value = "64.71"
$ws.Range("E17").Value = "  +0.57%  "
$ws.Range("D18").Value = "26.306.18"
$ws.Range("E18").Value = "  -1.02%  "
$ws.Range("D19").Value = "4.965"
$ws.Range("E19").Value = "  +0.59%  "
$ws.Range("E20").Value = "  -0.30%  "
$ws.Range("D21").Value = "10.86"
$ws.Range("E21").Value = "  +3.15%  "
$ws.Range("D22").Value = "189.65"
$ws.Range("E22").Value = "  +5.65%  "
$ws.Range("D23").Value = "6.221"
$ws.Range("E23").Value = "  +0.57%  "
$ws.Range("E24").Value = "  -0.29%  "
$ws.Range("D25").Value = "144.56"
$ws.Range("E25").Value = "  +0.98%  "
$ws.Range("D26").Value = "7.640"
$ws.Range("E26").Value = "  +1.09%  "
$ws.Range("D27").Value = "0.1187"
$ws.Range("E27").Value = "  +3.92%  "
$ws.Range("D28").Value = "15.64"
$ws.Range("E28").Value = "  +4.71%  "
$ws.Range("D29").Value = "0.05974"
$ws.Range("E29").Value = "  +2.96%  "
$ws.Range("D30").Value = "1.284"
$ws.Range("E30").Value = "  -3.46%  "
$ws.Range("D31").Value = "1.320"
$ws.Range("E31").Value = "  -1.78%  "
$ws.Range("D32").Value = "3.530"
$ws.Range("E32").Value = "  +3.48%  "
$ws.Range("D33").Value = "3.530"
$ws.Range("E33").Value = "  +3.94%  "
$ws.Range("D34").Value = "1.638"
$ws.Range("E34").Value = "  +4.04%  "
$ws.Range("D35").Value = "1.016"
$ws.Range("E35").Value = "  +4.46%  "
$ws.Range("D36").Value = "0.6028"
$ws.Range("E36").Value = "  +1.79%  "
$ws.Range("D37").Value = "2.374"
$ws.Range("E37").Value = "  -2.31%  "
$ws.Range("D38").Value = "2.648"
$ws.Range("E38").Value = "  +0.68%  "
$ws.Range("D39").Value = "0.01621"
$ws.Range("E39").Value = "  +3.07%  "
$ws.Range("D40").Value = "6.081"
$ws.Range("E40").Value = "  +5.98%  "
$ws.Range("D41").Value = "1.080.47"
$ws.Range("E41").Value = "  +0.97%  "
$ws.Range("D42").Value = "0.8707"
$ws.Range("E42").Value = "  +1.00%  "
$ws.Range("D43").Value = "1.011"
$ws.Range("E43").Value = "  -0.09%  "
$ws.Range("D44").Value = "100.17"
$ws.Range("E44").Value = "  +4.02%  "
$ws.Range("D45").Value = "1.821.59"
$ws.Range("E45").Value = "  +1.51%  "
$ws.Range("D46").Value = "0.00000000115"
$ws.Range("E46").Value = "  +9.36%  "
$ws.Range("D47").Value = "56.18"
$ws.Range("E47").Value = "  +1.99%  "
$ws.Range("D48").Value = "1.008"
$ws.Range("E48").Value = "  +0.19%  "
$ws.Range("D49").Value = "8.046"
$ws.Range("E49").Value = "  +3.56%  "
$ws.Range("D50").Value = "0.05206"
$ws.Range("E50").Value = "  +0.19%  "
$ws.Range("D51").Value = "0.4299"
$ws.Range("E51").Value = "  -1.85%  "

$ws.Range("D2:D51").ClearFormats()
